$d = $word.ActiveDocument

# The first paragraph used to read:
#   "Missing pipe arter variable name :    <---M2Doc version mismatch: template
#    is 3.1.1 and runtime is 3.2.0    "
# spread across several runs (plain spacer runs plus a highlighted "<---" /
# message pair). The template version mismatch has been fixed, so drop
# everything that follows the leading sentence, keeping only the first run:
#   "Missing pipe arter variable name :"
$keep = "Missing pipe arter variable name :"
$para = $d.Paragraphs(1).Range

$start = $para.Start + $keep.Length
$end = $para.End - 1          # stop before the paragraph mark

if ($end -gt $start) {
    $extra = $d.Range($start, $end)
    $extra.Delete()
}
